$d = $word.ActiveDocument

$pairs = @(
    @("959÷9=", "261÷3="),
    @("650÷8=", "618÷6="),
    @("748÷4=", "811÷7="),
    @("939÷6=", "436÷6="),
    @("752÷5=", "796÷6="),
    @("461÷6=", "746÷9="),
    @("180÷3=", "113÷4="),
    @("155÷7=", "351÷8="),
    @("981÷6=", "284÷4="),
    @("914÷2=", "904÷9="),
    @("965÷8=", "328÷8="),
    @("346÷4=", "199÷5="),
    @("580÷9=", "749÷9="),
    @("609÷6=", "301÷3="),
    @("148÷4=", "160÷2="),
    @("514÷2=", "851÷7="),
    @("627÷8=", "511÷7="),
    @("686÷9=", "276÷5="),
    @("338÷8=", "706÷6="),
    @("690÷4=", "955÷6="),
    @("761÷3=", "274÷3="),
    @("718÷8=", "631÷8="),
    @("625÷5=", "476÷2="),
    @("116÷3=", "460÷6="),
    @("119÷2=", "814÷2=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
